$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stundenerfassung")

# --- Row 6: hours corrected 2.5 -> 4.5 ---
$ws.Range("D6").Value = 4.5

# --- New rows 7-15: dates/employee copy formatting from an existing date cell ---
$ws.Range("A2").Copy($ws.Range("A7:A9"))
$ws.Range("A2").Copy($ws.Range("A10"))
$ws.Range("A2").Copy($ws.Range("A11:A15"))

$ws.Range("A7").Value = 43686
$ws.Range("A8").Value = 43689
$ws.Range("A9").Value = 43693
$ws.Range("A10").Value = 43714
$ws.Range("A11").Value = 43715
$ws.Range("A12").Value = 43721
$ws.Range("A13").Value = 43723
$ws.Range("A14").Value = 43732
$ws.Range("A15").Value = 43734

$ws.Range("B7:B15").Value = "MG"

# Work-description strings must be entered in this exact order so the
# shared-string table indices line up with the author's original edit
# sequence (row 10 / "ERM" was filled in after rows 11-14).
$ws.Range("C7").Value = "MVC Konzpt"
$ws.Range("C8").Value = "Webpack aufsetzen"
$ws.Range("C9").Value = "Dispatching"
$ws.Range("C11").Value = "Umsetzung Datenbankanbindung"
$ws.Range("C12").Value = "Umsetzung Site and Date Managment "
$ws.Range("C13").Value = "Umsetzung Comany und User Managment"
$ws.Range("C14").Value = "Umsetzung Reservation Managment"
$ws.Range("C10").Value = "ERM"
$ws.Range("C15").Value = "Präsentation"

$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 3
$ws.Range("D9").Value = 5
$ws.Range("D10").Value = 3
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 6
$ws.Range("D13").Value = 6
$ws.Range("D14").Value = 6
$ws.Range("D15").Value = 4

# --- Totals row (row 16 stays empty, total lands on row 17) ---
$ws.Range("D17").Formula = "=SUM(D2:D15)"

# --- Selection moved to H18 as left by the author on save ---
$ws.Range("H18").Select()
